$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "subject" header in column E, matching the bold style of the other headers
$ws.Range("E1").Value = "subject"
$ws.Range("E1").Font.Bold = $true

# Add formulas that build the subject line from the property column (C)
$ws.Range("E2").Formula = "=CONCAT(C2, "" - Owner Statement"")"
$ws.Range("E3").Formula = "=CONCAT(C3, "" - Owner Statement"")"
$ws.Range("E4").Formula = "=CONCAT(C4, "" - Owner Statement"")"
$ws.Range("E5").Formula = "=CONCAT(C5, "" - Owner Statement"")"

# Set the width of the new column E (closest achievable value to the
# target 25.90625 character-width units, given this engine's column-width
# rounding granularity)
$ws.Columns("E").ColumnWidth = 25

# Move the active selection to A8, matching the updated sheet view
$ws.Range("A8").Select()
